$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.209.25"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.24"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.60"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6079"
$ws.Range("E6").Value = "  -3.44%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07102"
$ws.Range("E8").Value = "  -4.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2812"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.77"
$ws.Range("E10").Value = "  -4.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.827.89"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.825"
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001011"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6356"
$ws.Range("E15").Value = "  -5.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.067.22"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.22"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.892"
$ws.Range("E18").Value = "  -5.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.200.50"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.92"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.80"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.019"
$ws.Range("E23").Value = "  -4.72%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.61"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.073"
$ws.Range("E26").Value = "  -4.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1302"
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.57"
$ws.Range("E28").Value = "  -4.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.479"
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06459"
$ws.Range("E30").Value = "  -6.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.825"
$ws.Range("E32").Value = "  -5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.808"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.128"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.754"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6507"
$ws.Range("E36").Value = "  -6.79%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.218.59"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01752"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.531"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9313"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.08"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.983.04"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "63.08"
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.613"
$ws.Range("E48").Value = "  -5.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.575"
$ws.Range("E49").Value = "  -4.07%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1076"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05527"
$ws.Range("E51").Value = "  -2.65%  "

Write-Host "Applied cryptos update"
